$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record needs to be inserted above the existing
# history (which starts at row 160), pushing every subsequent record
# down by two rows (one "Primera" + one "Segunda" row per date).
$ws.Rows.Item(160).Insert()
$ws.Rows.Item(160).Insert()

# The two freshly inserted rows are blank (format-only). Seed them with
# the record that is now sitting right below them (i.e. the record that
# used to occupy rows 160-161 before the insert), then overwrite the
# fields that actually differ for the new week.
$ws.Range("A162:R162").Copy()
$ws.Range("A160").PasteSpecial()
$ws.Range("A163:R163").Copy()
$ws.Range("A161").PasteSpecial()

# New record's own date / price.
$ws.Range("D160").Value = 44466
$ws.Range("D161").Value = 44466
$ws.Range("J161").Value = 1800
